$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.607.83'
$ws.Range("E2").Value = '  +2.29%  '
$ws.Range("D3").Value = '3.057.68'
$ws.Range("E3").Value = '  +2.57%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''518.42'
$ws.Range("E5").Value = '  +2.67%  '
$ws.Range("D6").Value = '''141.76'
$ws.Range("E6").Value = '  +3.77%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  +1.32%  '
$ws.Range("D9").Value = '''7.28'
$ws.Range("E9").Value = '  +1.31%  '
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("E11").Value = '  +3.02%  '
$ws.Range("D12").Value = '3.576.41'
$ws.Range("E12").Value = '  +2.37%  '
$ws.Range("D13").Value = '''0.130'
$ws.Range("E13").Value = '  +3.24%  '
$ws.Range("D14").Value = '''25.60'
$ws.Range("E14").Value = '  -0.80%  '
$ws.Range("D15").Value = '''0.0000163'
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").Value = '57.648.85'
$ws.Range("E16").Value = '  +2.46%  '
$ws.Range("D17").Value = '3.050.53'
$ws.Range("E17").Value = '  +2.26%  '
$ws.Range("D18").Value = '''6.07'
$ws.Range("E18").Value = '  +1.62%  '
$ws.Range("D19").Value = '''12.79'
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("D20").Value = '''8.14'
$ws.Range("E20").Value = '  +1.50%  '
$ws.Range("D21").Value = '''330.09'
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = '''0.496'
$ws.Range("E23").Value = '  +0.89%  '
$ws.Range("D24").Value = '''65.69'
$ws.Range("E24").Value = '  +1.92%  '
$ws.Range("E25").Value = '  +3.83%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = '0.0₃0898'
$ws.Range("E27").Value = '  -2.07%  '
$ws.Range("D28").Value = '''6.33'
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("D29").Value = '''7.19'
$ws.Range("E29").Value = '  +3.50%  '
$ws.Range("E30").Value = '  +2.38%  '
$ws.Range("D31").Value = '''20.70'
$ws.Range("E31").Value = '  +2.81%  '
$ws.Range("E32").Value = '  +2.52%  '
$ws.Range("D33").Value = '''154.30'
$ws.Range("E33").Value = '  +1.06%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = '''4.48'
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("B35").Value = 'EnergySwap'
$ws.Range("C35").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D35").Value = '''27.17'
$ws.Range("E35").Value = '  +4.92%  '
$ws.Range("D36").Value = '''5.92'
$ws.Range("E36").Value = '  +2.78%  '
$ws.Range("D37").Value = '''1.25'
$ws.Range("E37").Value = '  +1.34%  '
$ws.Range("D38").Value = '''0.0672'
$ws.Range("E38").Value = '  +2.12%  '
$ws.Range("D39").Value = '3.091.64'
$ws.Range("E39").Value = '  +2.41%  '
$ws.Range("D40").Value = '''3.91'
$ws.Range("E40").Value = '  +3.42%  '
$ws.Range("D41").Value = '''36.78'
$ws.Range("E41").Value = '  -0.41%  '
$ws.Range("D42").Value = '''0.999'
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("D44").Value = '2.254.97'
$ws.Range("E44").Value = '  +3.91%  '
$ws.Range("D45").Value = '''0.0256'
$ws.Range("E45").Value = '  +9.09%  '
$ws.Range("D46").Value = '''20.72'
$ws.Range("E46").Value = '  +6.83%  '
$ws.Range("D47").Value = '''1.36'
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").Value = '''5.86'
$ws.Range("E48").Value = '  +0.89%  '
$ws.Range("D49").Value = '''0.916'
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("D50").Value = '''262.25'
$ws.Range("E50").Value = '  +15.62%  '
$ws.Range("D51").Value = '''0.715'
$ws.Range("E51").Value = '  +6.59%  '
